# Auto-generated edit script: updates cached market-price columns (H:N)
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets to match refreshed API data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value2 = 4000
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 4000
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 4000
$ws.Range("M64").Value2 = ""
$ws.Range("N64").Value2 = -4496

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value2 = 4000
$ws.Range("I67").Value2 = 0
$ws.Range("J67").Value2 = 4000
$ws.Range("K67").Value2 = 0
$ws.Range("L67").Value2 = 4000
$ws.Range("M67").Value2 = ""
$ws.Range("N67").Value2 = -5716

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value2 = 905.94116
$ws.Range("I129").Value2 = 1286.3334
$ws.Range("J129").Value2 = 882.1667
$ws.Range("K129").Value2 = 3859.0002
$ws.Range("L129").Value2 = 2646.5001
$ws.Range("M129").Value2 = 1140.9998
$ws.Range("N129").Value2 = -12646.5001

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value2 = 1243.9333
$ws.Range("I137").Value2 = 1189.1538
$ws.Range("J137").Value2 = 1600
$ws.Range("K137").Value2 = 3567.4614
$ws.Range("L137").Value2 = 4800
$ws.Range("M137").Value2 = -1017.4614
$ws.Range("N137").Value2 = -9900

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value2 = 2789.4902
$ws.Range("I138").Value2 = 2685.923
$ws.Range("K138").Value2 = 8057.768999999999
$ws.Range("M138").Value2 = -2917.768999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value2 = 3511.0728
$ws.Range("I32").Value2 = 2497.8086
$ws.Range("J32").Value2 = 9464
$ws.Range("K32").Value2 = 2497.8086
$ws.Range("L32").Value2 = 9464
$ws.Range("M32").Value2 = -2210.8086
$ws.Range("N32").Value2 = -10038

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value2 = 3373.2068
$ws.Range("I61").Value2 = 2708.6667
$ws.Range("K61").Value2 = 2708.6667
$ws.Range("M61").Value2 = -2496.6667

# Row 88 (Leve Item ID 12530)
$ws.Range("I88").Value2 = 1999
$ws.Range("J88").Value2 = 4642.5713
$ws.Range("K88").Value2 = 1999
$ws.Range("L88").Value2 = 4642.5713
$ws.Range("M88").Value2 = -1593
$ws.Range("N88").Value2 = -5454.5713

# Row 91 (Leve Item ID 12530)
$ws.Range("I91").Value2 = 1999
$ws.Range("J91").Value2 = 4642.5713
$ws.Range("K91").Value2 = 1999
$ws.Range("L91").Value2 = 4642.5713
$ws.Range("M91").Value2 = -595
$ws.Range("N91").Value2 = -7450.5713

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value2 = 1125.8
$ws.Range("I97").Value2 = 1125.8
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 1125.8
$ws.Range("L97").Value2 = 0
$ws.Range("M97").Value2 = -629.8
$ws.Range("N97").Value2 = ""

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value2 = 2503.4285
$ws.Range("J102").Value2 = 3500
$ws.Range("L102").Value2 = 3500
$ws.Range("N102").Value2 = -6744

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value2 = 1552.6471
$ws.Range("I132").Value2 = 951.7895
$ws.Range("J132").Value2 = 2313.7334
$ws.Range("K132").Value2 = 2855.3685
$ws.Range("L132").Value2 = 6941.2002
$ws.Range("M132").Value2 = -325.3685
$ws.Range("N132").Value2 = -12001.2002

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value2 = 3373.2068
$ws.Range("I136").Value2 = 2708.6667
$ws.Range("K136").Value2 = 8126.000100000001
$ws.Range("M136").Value2 = -5576.000100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value2 = 1854.9546
$ws.Range("I20").Value2 = 1877.4667
$ws.Range("J20").Value2 = 1806.7142
$ws.Range("K20").Value2 = 1877.4667
$ws.Range("L20").Value2 = 1806.7142
$ws.Range("M20").Value2 = -1630.4667
$ws.Range("N20").Value2 = -2300.7142

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value2 = 97360.81
$ws.Range("I86").Value2 = 1938.5333
$ws.Range("K86").Value2 = 1938.5333
$ws.Range("M86").Value2 = -815.5333000000001

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value2 = 97360.81
$ws.Range("I89").Value2 = 1938.5333
$ws.Range("K89").Value2 = 9692.666499999999
$ws.Range("M89").Value2 = -4076.666499999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value2 = 2886.257
$ws.Range("I31").Value2 = 1952.1818
$ws.Range("J31").Value2 = 4467
$ws.Range("K31").Value2 = 1952.1818
$ws.Range("L31").Value2 = 4467
$ws.Range("M31").Value2 = -1657.1818
$ws.Range("N31").Value2 = -5057

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value2 = 2886.257
$ws.Range("I34").Value2 = 1952.1818
$ws.Range("J34").Value2 = 4467
$ws.Range("K34").Value2 = 1952.1818
$ws.Range("L34").Value2 = 4467
$ws.Range("M34").Value2 = -1750.1818
$ws.Range("N34").Value2 = -4871

# Row 43 (Leve Item ID 18504)
$ws.Range("H43").Value2 = 16249.75
$ws.Range("J43").Value2 = 16249.75
$ws.Range("L43").Value2 = 16249.75
$ws.Range("N43").Value2 = -16617.75

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value2 = 1673721.6
$ws.Range("I58").Value2 = 3106640.8
$ws.Range("J58").Value2 = 1982.8334
$ws.Range("K58").Value2 = 3106640.8
$ws.Range("L58").Value2 = 1982.8334
$ws.Range("M58").Value2 = -3106437.8
$ws.Range("N58").Value2 = -2388.8334

# Row 101 (Leve Item ID 18504)
$ws.Range("H101").Value2 = 16249.75
$ws.Range("J101").Value2 = 16249.75
$ws.Range("L101").Value2 = 16249.75
$ws.Range("N101").Value2 = -22739.75

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value2 = 1673721.6
$ws.Range("I136").Value2 = 3106640.8
$ws.Range("J136").Value2 = 1982.8334
$ws.Range("K136").Value2 = 9319922.399999999
$ws.Range("L136").Value2 = 5948.5002
$ws.Range("M136").Value2 = -9317372.399999999
$ws.Range("N136").Value2 = -11048.5002

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value2 = 3323.5
$ws.Range("I80").Value2 = 3397
$ws.Range("K80").Value2 = 3397
$ws.Range("M80").Value2 = -2399

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value2 = 3323.5
$ws.Range("I83").Value2 = 3397
$ws.Range("K83").Value2 = 16985
$ws.Range("M83").Value2 = -11993

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value2 = 1654.9375
$ws.Range("I122").Value2 = 1312.1
$ws.Range("K122").Value2 = 3936.3
$ws.Range("M122").Value2 = -1486.3

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value2 = 1573410.2
$ws.Range("I126").Value2 = 2139339.5
$ws.Range("J126").Value2 = 101994.3
$ws.Range("K126").Value2 = 6418018.5
$ws.Range("L126").Value2 = 305982.9
$ws.Range("M126").Value2 = -6415548.5
$ws.Range("N126").Value2 = -310922.9

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value2 = 1168202.1
$ws.Range("I132").Value2 = 1426641.9
$ws.Range("J132").Value2 = 5223.3335
$ws.Range("K132").Value2 = 4279925.699999999
$ws.Range("L132").Value2 = 15670.0005
$ws.Range("M132").Value2 = -4277395.699999999
$ws.Range("N132").Value2 = -20730.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (Leve Item ID 2631)
$ws.Range("H2").Value2 = 418000
$ws.Range("J2").Value2 = 90000
$ws.Range("L2").Value2 = 90000
$ws.Range("N2").Value2 = -90224

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value2 = 2327.1428
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 2327.1428
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 2327.1428
$ws.Range("M46").Value2 = ""
$ws.Range("N46").Value2 = -2703.1428

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value2 = 286.9697
$ws.Range("J55").Value2 = 356.2
$ws.Range("L55").Value2 = 356.2
$ws.Range("N55").Value2 = -702.2

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value2 = 3716.6875
$ws.Range("I136").Value2 = 1744.8334
$ws.Range("J136").Value2 = 4899.8
$ws.Range("K136").Value2 = 5234.5002
$ws.Range("L136").Value2 = 14699.4
$ws.Range("M136").Value2 = -2684.5002
$ws.Range("N136").Value2 = -19799.4

# Row 139 (Leve Item ID 43310)
$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").Value2 = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 40 (Leve Item ID 3601)
$ws.Range("H40").Value2 = 0
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 0
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = ""
$ws.Range("N40").Value2 = ""

# Row 63 (Leve Item ID 10824)
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("N63").Value2 = ""

# Row 66 (Leve Item ID 10824)
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("N66").Value2 = ""

# Row 69 (Leve Item ID 10951)
$ws.Range("H69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("N69").Value2 = ""

# Row 72 (Leve Item ID 10951)
$ws.Range("H72").Value2 = 0
$ws.Range("J72").Value2 = 0
$ws.Range("L72").Value2 = 0
$ws.Range("N72").Value2 = ""

# Row 75 (Leve Item ID 11957)
$ws.Range("H75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("N75").Value2 = ""

# Row 76 (Leve Item ID 10896)
$ws.Range("H76").Value2 = 0
$ws.Range("J76").Value2 = 0
$ws.Range("L76").Value2 = 0
$ws.Range("N76").Value2 = ""

# Row 78 (Leve Item ID 11957)
$ws.Range("H78").Value2 = 0
$ws.Range("J78").Value2 = 0
$ws.Range("L78").Value2 = 0
$ws.Range("N78").Value2 = ""

# Row 79 (Leve Item ID 10896)
$ws.Range("H79").Value2 = 0
$ws.Range("J79").Value2 = 0
$ws.Range("L79").Value2 = 0
$ws.Range("N79").Value2 = ""

# Row 80 (Leve Item ID 10911)
$ws.Range("H80").Value2 = 79998
$ws.Range("J80").Value2 = 79998
$ws.Range("L80").Value2 = 79998
$ws.Range("N80").Value2 = -81994

# Row 83 (Leve Item ID 10911)
$ws.Range("H83").Value2 = 79998
$ws.Range("J83").Value2 = 79998
$ws.Range("L83").Value2 = 239994
$ws.Range("N83").Value2 = -249978

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value2 = 22500.428
$ws.Range("I126").Value2 = 27500.6
$ws.Range("J126").Value2 = 10000
$ws.Range("K126").Value2 = 82501.79999999999
$ws.Range("L126").Value2 = 30000
$ws.Range("M126").Value2 = -80031.79999999999
$ws.Range("N126").Value2 = -34940

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value2 = 2776.64
$ws.Range("I132").Value2 = 2411.2856
$ws.Range("K132").Value2 = 7233.8568
$ws.Range("M132").Value2 = -4703.8568

# Row 139 (Leve Item ID 43312)
$ws.Range("H139").Value2 = 59982.5
$ws.Range("I139").Value2 = 0
$ws.Range("J139").Value2 = 59982.5
$ws.Range("K139").Value2 = 0
$ws.Range("L139").Value2 = 59982.5
$ws.Range("M139").Value2 = ""
$ws.Range("N139").Value2 = -70262.5
